$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 19 (Leve Item ID 7015)
$ws.Range("H19").Value = 2291.5
$ws.Range("I19").Value = 2283.5
$ws.Range("K19").Value = 2283.5
$ws.Range("M19").Value = -2108.5
# row 28 (Leve Item ID 27772)
$ws.Range("H28").Value = 2492.9
$ws.Range("I28").Value = 2362.2856
$ws.Range("J28").Value = 2797.6667
$ws.Range("K28").Value = 2362.2856
$ws.Range("L28").Value = 2797.6667
$ws.Range("M28").Value = -1877.2856
$ws.Range("N28").Value = -3767.6667
# row 32 (Leve Item ID 5484)
$ws.Range("H32").Value = 4563.4
$ws.Range("J32").Value = 4563.4
$ws.Range("L32").Value = 4563.4
$ws.Range("N32").Value = -5215.4
# row 55 (Leve Item ID 5517)
$ws.Range("H55").Value = 768.25
$ws.Range("I55").Value = 814.6667
$ws.Range("J55").Value = 740.4
$ws.Range("K55").Value = 814.6667
$ws.Range("L55").Value = 740.4
$ws.Range("M55").Value = -600.6667
$ws.Range("N55").Value = -1168.4
# row 92 (Leve Item ID 19901)
$ws.Range("H92").Value = 848.6667
$ws.Range("I92").Value = 810.35297
$ws.Range("J92").Value = 1500
$ws.Range("K92").Value = 810.35297
$ws.Range("L92").Value = 1500
$ws.Range("M92").Value = 437.64703
$ws.Range("N92").Value = -3996
# row 100 (Leve Item ID 19906)
$ws.Range("H100").Value = 1449.6666
$ws.Range("I100").Value = 1339.6
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 1339.6
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -798.5999999999999
$ws.Range("N100").Value = -3082
# row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 24999
$ws.Range("I138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("M138").Value = $null

$ws = $wb.Worksheets.Item("ARM")
# row 2 (Leve Item ID 27713)
$ws.Range("H2").Value = 2521.6
$ws.Range("I2").Value = 2521.6
$ws.Range("K2").Value = 2521.6
$ws.Range("M2").Value = -2408.6
# row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 4194.0312
$ws.Range("I32").Value = 4194.0312
$ws.Range("K32").Value = 4194.0312
$ws.Range("M32").Value = -3907.0312
# row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 3295
$ws.Range("I61").Value = 2884.3333
$ws.Range("K61").Value = 2884.3333
$ws.Range("M61").Value = -2672.3333
# row 74 (Leve Item ID 44000)
$ws.Range("H74").Value = 1500
$ws.Range("I74").Value = 1500
$ws.Range("K74").Value = 1500
$ws.Range("M74").Value = -626
# row 77 (Leve Item ID 44000)
$ws.Range("H77").Value = 1500
$ws.Range("I77").Value = 1500
$ws.Range("K77").Value = 7500
$ws.Range("M77").Value = -3132
# row 116 (Leve Item ID 27713)
$ws.Range("H116").Value = 2521.6
$ws.Range("I116").Value = 2521.6
$ws.Range("K116").Value = 2521.6
$ws.Range("M116").Value = -227.5999999999999
# row 122 (Leve Item ID 36168)
$ws.Range("H122").Value = 1153.875
$ws.Range("I122").Value = 1153.875
$ws.Range("K122").Value = 3461.625
$ws.Range("M122").Value = -1011.625
# row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 3295
$ws.Range("I136").Value = 2884.3333
$ws.Range("K136").Value = 8652.999899999999
$ws.Range("M136").Value = -6102.999899999999
# row 139 (Leve Item ID 42321)
$ws.Range("H139").Value = 99949.5
$ws.Range("J139").Value = 99949.5
$ws.Range("L139").Value = 99949.5
$ws.Range("N139").Value = -110229.5

$ws = $wb.Worksheets.Item("BSM")
# row 3 (Leve Item ID 27713)
$ws.Range("H3").Value = 2521.6
$ws.Range("I3").Value = 2521.6
$ws.Range("K3").Value = 2521.6
$ws.Range("M3").Value = -2407.6
# row 64 (Leve Item ID 14184)
$ws.Range("H64").Value = 2176.1333
$ws.Range("I64").Value = 1059.8572
$ws.Range("J64").Value = 3152.875
$ws.Range("K64").Value = 1059.8572
$ws.Range("L64").Value = 3152.875
$ws.Range("M64").Value = -834.8571999999999
$ws.Range("N64").Value = -3602.875
# row 67 (Leve Item ID 14184)
$ws.Range("H67").Value = 2176.1333
$ws.Range("I67").Value = 1059.8572
$ws.Range("J67").Value = 3152.875
$ws.Range("K67").Value = 1059.8572
$ws.Range("L67").Value = 3152.875
$ws.Range("M67").Value = -279.8571999999999
$ws.Range("N67").Value = -4712.875
# row 86 (Leve Item ID 12526)
$ws.Range("H86").Value = 9625.5625
$ws.Range("I86").Value = 10287.077
$ws.Range("J86").Value = 6759
$ws.Range("K86").Value = 10287.077
$ws.Range("L86").Value = 6759
$ws.Range("M86").Value = -9164.076999999999
$ws.Range("N86").Value = -9005
# row 89 (Leve Item ID 12526)
$ws.Range("H89").Value = 9625.5625
$ws.Range("I89").Value = 10287.077
$ws.Range("J89").Value = 6759
$ws.Range("K89").Value = 51435.38499999999
$ws.Range("L89").Value = 33795
$ws.Range("M89").Value = -45819.38499999999
$ws.Range("N89").Value = -45027
# row 94 (Leve Item ID 19939)
$ws.Range("H94").Value = 842.6
$ws.Range("I94").Value = 786.6429000000001
$ws.Range("J94").Value = 891.5625
$ws.Range("K94").Value = 786.6429000000001
$ws.Range("L94").Value = 891.5625
$ws.Range("M94").Value = -335.6429000000001
$ws.Range("N94").Value = -1793.5625
# row 107 (Leve Item ID 27706)
$ws.Range("H107").Value = 2069.9092
$ws.Range("I107").Value = 2057.6
$ws.Range("J107").Value = 2193
$ws.Range("K107").Value = 2057.6
$ws.Range("L107").Value = 2193
$ws.Range("M107").Value = -137.5999999999999
$ws.Range("N107").Value = -6033

$ws = $wb.Worksheets.Item("CRP")
# row 140 (Leve Item ID 42455)
$ws.Range("H140").Value = 71899
$ws.Range("J140").Value = 71899
$ws.Range("L140").Value = 71899
$ws.Range("N140").Value = -82259

$ws = $wb.Worksheets.Item("CUL")
# row 5 (Leve Item ID 43974)
$ws.Range("H5").Value = 690.875
$ws.Range("I5").Value = 690.875
$ws.Range("K5").Value = 2072.625
$ws.Range("M5").Value = -1960.625
# row 92 (Leve Item ID 19841)
$ws.Range("H92").Value = 312.7143
$ws.Range("I92").Value = 312.7143
$ws.Range("K92").Value = 938.1428999999999
$ws.Range("M92").Value = 309.8571000000001
# row 111 (Leve Item ID 27856)
$ws.Range("H111").Value = 6174.75
$ws.Range("I111").Value = 6174.75
$ws.Range("K111").Value = 18524.25
$ws.Range("M111").Value = -15457.25
# row 116 (Leve Item ID 27866)
$ws.Range("H116").Value = 2000
$ws.Range("J116").Value = 2000
$ws.Range("L116").Value = 6000
$ws.Range("N116").Value = -12884
# row 131 (Leve Item ID 36060)
$ws.Range("H131").Value = 1000.8571
$ws.Range("I131").Value = 551.75
$ws.Range("J131").Value = 1599.6666
$ws.Range("K131").Value = 1655.25
$ws.Range("L131").Value = 4798.9998
$ws.Range("M131").Value = 3384.75
$ws.Range("N131").Value = -14878.9998
# row 132 (Leve Item ID 43972)
$ws.Range("H132").Value = 2166.3333
$ws.Range("I132").Value = 1999.5
$ws.Range("K132").Value = 17995.5
$ws.Range("M132").Value = -15465.5
# row 135 (Leve Item ID 43974)
$ws.Range("H135").Value = 690.875
$ws.Range("I135").Value = 690.875
$ws.Range("K135").Value = 6217.875
$ws.Range("M135").Value = -3682.875
# row 136 (Leve Item ID 44093)
$ws.Range("H136").Value = 6251
$ws.Range("I136").Value = 6251
$ws.Range("K136").Value = 18753
$ws.Range("M136").Value = -13653

$ws = $wb.Worksheets.Item("GSM")
# row 31 (Leve Item ID 2118)
$ws.Range("H31").Value = 1245
$ws.Range("I31").Value = 1245
$ws.Range("K31").Value = 1245
$ws.Range("M31").Value = -953
# row 37 (Leve Item ID 2118)
$ws.Range("H37").Value = 1245
$ws.Range("I37").Value = 1245
$ws.Range("K37").Value = 1245
$ws.Range("M37").Value = -968
# row 55 (Leve Item ID 4237)
$ws.Range("H55").Value = 2500
$ws.Range("I55").Value = 2500
$ws.Range("K55").Value = 2500
$ws.Range("M55").Value = -2173
# row 70 (Leve Item ID 14146)
$ws.Range("H70").Value = 3966.0908
$ws.Range("I70").Value = 5069.6
$ws.Range("K70").Value = 5069.6
$ws.Range("M70").Value = -4799.6
# row 73 (Leve Item ID 14146)
$ws.Range("H73").Value = 3966.0908
$ws.Range("I73").Value = 5069.6
$ws.Range("K73").Value = 5069.6
$ws.Range("M73").Value = -4133.6
# row 97 (Leve Item ID 19940)
$ws.Range("H97").Value = 1271.125
$ws.Range("I97").Value = 1294.8334
$ws.Range("J97").Value = 1200
$ws.Range("K97").Value = 1294.8334
$ws.Range("L97").Value = 1200
$ws.Range("M97").Value = -798.8334
$ws.Range("N97").Value = -2192
# row 113 (Leve Item ID 27710)
$ws.Range("H113").Value = 3293.4
$ws.Range("I113").Value = 3366.75
$ws.Range("K113").Value = 3366.75
$ws.Range("M113").Value = -1196.75
# row 122 (Leve Item ID 36182)
$ws.Range("H122").Value = 5055.6
$ws.Range("I122").Value = 4499.5
$ws.Range("J122").Value = 5194.625
$ws.Range("K122").Value = 13498.5
$ws.Range("L122").Value = 15583.875
$ws.Range("M122").Value = -11048.5
$ws.Range("N122").Value = -20483.875
# row 126 (Leve Item ID 36184)
$ws.Range("H126").Value = 3368.1
$ws.Range("I126").Value = 3368.1
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 10104.3
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -7634.299999999999
$ws.Range("N126").Value = $null
# row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 3571.5715
$ws.Range("I132").Value = 3514.5
$ws.Range("K132").Value = 10543.5
$ws.Range("M132").Value = -8013.5
# row 138 (Leve Item ID 42325)
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = $null

$ws = $wb.Worksheets.Item("LTW")
# row 22 (Leve Item ID 5277)
$ws.Range("H22").Value = 4079.875
$ws.Range("I22").Value = 3720
$ws.Range("J22").Value = 4439.75
$ws.Range("K22").Value = 3720
$ws.Range("L22").Value = 4439.75
$ws.Range("M22").Value = -3425
$ws.Range("N22").Value = -5029.75
# row 27 (Leve Item ID 5277)
$ws.Range("H27").Value = 4079.875
$ws.Range("I27").Value = 3720
$ws.Range("J27").Value = 4439.75
$ws.Range("K27").Value = 3720
$ws.Range("L27").Value = 4439.75
$ws.Range("M27").Value = -3613
$ws.Range("N27").Value = -4653.75
# row 61 (Leve Item ID 27740)
$ws.Range("H61").Value = 2402.0715
$ws.Range("I61").Value = 2371.7778
$ws.Range("K61").Value = 2371.7778
$ws.Range("M61").Value = -2169.7778
# row 113 (Leve Item ID 27740)
$ws.Range("H113").Value = 2402.0715
$ws.Range("I113").Value = 2371.7778
$ws.Range("K113").Value = 2371.7778
$ws.Range("M113").Value = -201.7777999999998
# row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 2984.9285
$ws.Range("I132").Value = 2945.3076
$ws.Range("K132").Value = 8835.9228
$ws.Range("M132").Value = -6305.9228
# row 136 (Leve Item ID 44060)
$ws.Range("H136").Value = 17724.818
$ws.Range("I136").Value = 17554.445
$ws.Range("K136").Value = 52663.335
$ws.Range("M136").Value = -50113.335

$ws = $wb.Worksheets.Item("WVR")
# row 113 (Leve Item ID 27752)
$ws.Range("H113").Value = 313.75
$ws.Range("I113").Value = 321.5
$ws.Range("K113").Value = 964.5
$ws.Range("M113").Value = 1205.5
# row 126 (Leve Item ID 36210)
$ws.Range("H126").Value = 2879.9092
$ws.Range("I126").Value = 4535.8
$ws.Range("K126").Value = 13607.4
$ws.Range("M126").Value = -11137.4
# row 136 (Leve Item ID 44031)
$ws.Range("H136").Value = 11072.583
$ws.Range("I136").Value = 12767.4
$ws.Range("K136").Value = 38302.2
$ws.Range("M136").Value = -35752.2
